$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Real Data")
$ws.Range("A1").Value = "Test"
